# Apply custom accuracy (round to 2 decimals) on row 5 and delete row 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, rounded (2 decimal place) values for B5:AH5
$row5Values = @(6.04, 4.88, 0.37, 13.24, 11.09, 4.93, 23.48, 7.51, 3.42, 5.43, 6.12, 5.62, 1.42, 4.66, 7.05, 4.03, 0.04, 0.15, 66.92, 13.82, 4.22, 9.24, 5.16, 0.66, 11.25, 3.89, 4.01, 3.93, 6.08, 0.08, 21.51, 2.66, 5.5)

$col = 2
foreach ($val in $row5Values) {
    $ws.Cells.Item(5, $col).Value = $val
    $col = $col + 1
}

# Remove the last data row (row 6), shifting rows up (there is nothing below it).
$ws.Rows.Item(6).Delete()
